$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.194828579899138
$ws.Range("C2").Value = 0.2170881049814568
$ws.Range("D2").Value = 0.02325209133015704
$ws.Range("E2").Value = 0.09891205464731279
$ws.Range("F2").Value = 0.7763925801715814
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("L2").Value = 0.212214618074313
$ws.Range("N2").Value = 1.152311375679432
$ws.Range("O2").Value = 2.66100384504972
$ws.Range("B3").Value = 1.090023758153961
$ws.Range("C3").Value = 0.2070080621694927
$ws.Range("D3").Value = 0.02224389641261837
$ws.Range("E3").Value = 0.09950659437024001
$ws.Range("F3").Value = 0.7690441028765775
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("L3").Value = 0.2029099639223659
$ws.Range("N3").Value = 1.16370447432535
$ws.Range("O3").Value = 2.651355265332029
$ws.Range("B4").Value = 1.025899365142493
$ws.Range("C4").Value = 0.2007727506100849
$ws.Range("D4").Value = 0.02162031459644354
$ws.Range("E4").Value = 0.09992322555521049
$ws.Range("F4").Value = 0.7650888328931629
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("L4").Value = 0.1973047040675198
$ws.Range("N4").Value = 1.171190191582255
$ws.Range("O4").Value = 2.647313203223092
$ws.Range("B5").Value = 0.9998263544598558
$ws.Range("C5").Value = 0.1982203687496309
$ws.Range("D5").Value = 0.0213650737288269
$ws.Range("E5").Value = 0.1001059769190888
$ws.Range("F5").Value = 0.7636168359997342
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("L5").Value = 0.1950476885298542
$ws.Range("N5").Value = 1.17436401216753
$ws.Range("O5").Value = 2.646138631301056
$ws.Range("B6").Value = 0.9955005050966292
$ws.Range("C6").Value = 0.1977958612505404
$ws.Range("D6").Value = 0.0213226236564239
$ws.Range("E6").Value = 0.1001371060547349
$ws.Range("F6").Value = 0.7633808506751265
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("L6").Value = 0.1946745557265643
$ws.Range("N6").Value = 1.174898472694672
$ws.Range("O6").Value = 2.64597211820103
$ws.Range("B7").Value = 1.025547497997707
$ws.Range("C7").Value = 0.2007383744144278
$ws.Range("D7").Value = 0.02161687686643177
$ws.Range("E7").Value = 0.09992563768537188
$ws.Range("F7").Value = 0.7650684151992451
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("L7").Value = 0.1972741550656139
$ws.Range("N7").Value = 1.171232495426999
$ws.Range("O7").Value = 2.647295449918147
$ws.Range("B8").Value = 1.158645690242849
$ws.Range("C8").Value = 0.2136221702924956
$ws.Range("D8").Value = 0.02290541910641508
$ws.Range("E8").Value = 0.09910634624981185
$ws.Range("F8").Value = 0.7737431608608674
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("L8").Value = 0.2089840230086253
$ws.Range("N8").Value = 1.156137974602188
$ws.Range("O8").Value = 2.657285933822692
$ws.Range("B9").Value = 1.421400947532618
$ws.Range("C9").Value = 0.2385155255672942
$ws.Range("D9").Value = 0.02539554518413212
$ws.Range("E9").Value = 0.09790903914424831
$ws.Range("F9").Value = 0.7951821495015423
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("L9").Value = 0.2328015675969368
$ws.Range("N9").Value = 1.130426422940822
$ws.Range("O9").Value = 2.691848309765248
$ws.Range("B10").Value = 1.615474074365636
$ws.Range("C10").Value = 0.2565719894211043
$ws.Range("D10").Value = 0.02720198638427718
$ws.Range("E10").Value = 0.09727905442996843
$ws.Range("F10").Value = 0.8136502656608684
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("L10").Value = 0.2508218888156648
$ws.Range("N10").Value = 1.113904723711698
$ws.Range("O10").Value = 2.726426362757309
$ws.Range("B11").Value = 1.703978700736286
$ws.Range("C11").Value = 0.2647346361290204
$ws.Range("D11").Value = 0.0280186437467691
$ws.Range("E11").Value = 0.09704671115884622
$ws.Range("F11").Value = 0.8226457974670183
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("L11").Value = 0.2591333208483917
$ws.Range("N11").Value = 1.106902547696258
$ws.Range("O11").Value = 2.744164293284825
$ws.Range("B12").Value = 1.737523650117055
$ws.Range("C12").Value = 0.2678181003560098
$ws.Range("D12").Value = 0.02832714209637999
$ws.Range("E12").Value = 0.09696653071968875
$ws.Range("F12").Value = 0.826137886513294
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("L12").Value = 0.2622969971709921
$ws.Range("N12").Value = 1.104324859715156
$ws.Range("O12").Value = 2.751170867853546
$ws.Range("B13").Value = 1.730297826171238
$ws.Range("C13").Value = 0.2671543598127073
$ws.Range("D13").Value = 0.0282607351990336
$ws.Range("E13").Value = 0.09698345192269286
$ws.Range("F13").Value = 0.8253819887171119
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("L13").Value = 0.2616149181192924
$ws.Range("N13").Value = 1.104876724474579
$ws.Range("O13").Value = 2.749648983206214
$ws.Range("B14").Value = 1.706737867512345
$ws.Range("C14").Value = 0.2649884669413041
$ws.Range("D14").Value = 0.0280440392693535
$ws.Range("E14").Value = 0.09703995826488843
$ws.Range("F14").Value = 0.8229313754082597
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("L14").Value = 0.2593932717276175
$ws.Range("N14").Value = 1.10668899884606
$ws.Range("O14").Value = 2.744734919531027
$ws.Range("B15").Value = 1.692310598760628
$ws.Range("C15").Value = 0.2636608063158405
$ws.Range("D15").Value = 0.02791120826901761
$ws.Range("E15").Value = 0.0970755862986401
$ws.Range("F15").Value = 0.821441468037662
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("L15").Value = 0.2580345723620923
$ws.Range("N15").Value = 1.107808691947341
$ws.Range("O15").Value = 2.741762653354385
$ws.Range("B16").Value = 1.609694388645153
$ws.Range("C16").Value = 0.2560374936823848
$ws.Range("D16").Value = 0.0271485116659278
$ws.Range("E16").Value = 0.09729533033498861
$ws.Range("F16").Value = 0.8130743638329108
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("L16").Value = 0.2502810038860446
$ws.Range("N16").Value = 1.114372669914133
$ws.Range("O16").Value = 2.725307628259515
$ws.Range("B17").Value = 1.559067288818142
$ws.Range("C17").Value = 0.2513475691720544
$ws.Range("D17").Value = 0.02667930216745873
$ws.Range("E17").Value = 0.09744403065693241
$ws.Range("F17").Value = 0.8080937843583627
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("L17").Value = 0.2455535692936337
$ws.Range("N17").Value = 1.118531026435591
$ws.Range("O17").Value = 2.715727941592206
$ws.Range("B18").Value = 1.529968693121248
$ws.Range("C18").Value = 0.248645227518864
$ws.Range("D18").Value = 0.02640894616511957
$ws.Range("E18").Value = 0.09753466427133795
$ws.Range("F18").Value = 0.8052850114918186
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("L18").Value = 0.2428451981714517
$ws.Range("N18").Value = 1.120971146377784
$ws.Range("O18").Value = 2.710406931092308
$ws.Range("B19").Value = 1.520120020252477
$ws.Range("C19").Value = 0.2477294382818229
$ws.Range("D19").Value = 0.02631732654386809
$ws.Range("E19").Value = 0.09756622796606251
$ws.Range("F19").Value = 0.8043436077378345
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("L19").Value = 0.2419300339969652
$ws.Range("N19").Value = 1.121805631386472
$ws.Range("O19").Value = 2.708637759218277
$ws.Range("B20").Value = 1.564454490507444
$ws.Range("C20").Value = 0.2518473198586264
$ws.Range("D20").Value = 0.02672930003939911
$ws.Range("E20").Value = 0.09742767286917875
$ws.Range("F20").Value = 0.8086181859662105
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("L20").Value = 0.2460557033023605
$ws.Range("N20").Value = 1.118083358771003
$ws.Range("O20").Value = 2.716728151697907
$ws.Range("B21").Value = 1.713657191149082
$ws.Range("C21").Value = 0.2656248485309618
$ws.Range("D21").Value = 0.02810770867469614
$ws.Range("E21").Value = 0.09702314918731503
$ws.Range("F21").Value = 0.8236488529773851
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("L21").Value = 0.2600453807628185
$ws.Range("N21").Value = 1.106154684515104
$ws.Range("O21").Value = 2.746170432752024
$ws.Range("B22").Value = 1.811344960305803
$ws.Range("C22").Value = 0.2745851564138206
$ws.Range("D22").Value = 0.02900418625663548
$ws.Range("E22").Value = 0.09680425247719526
$ws.Range("F22").Value = 0.8339717507095514
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("L22").Value = 0.2692835256879391
$ws.Range("N22").Value = 1.098789274595937
$ws.Range("O22").Value = 2.767101089965735
$ws.Range("B23").Value = 1.759191604286855
$ws.Range("C23").Value = 0.2698069635567606
$ws.Range("D23").Value = 0.0285261272538051
$ws.Range("E23").Value = 0.09691691899231536
$ws.Range("F23").Value = 0.8284164510124015
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("L23").Value = 0.2643442739981481
$ws.Range("N23").Value = 1.102680916745427
$ws.Range("O23").Value = 2.755775239953834
$ws.Range("B24").Value = 1.562018912846838
$ws.Range("C24").Value = 0.2516214013805609
$ws.Range("D24").Value = 0.02670669787237756
$ws.Range("E24").Value = 0.09743505220290416
$ws.Range("F24").Value = 0.8083809338497474
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("L24").Value = 0.2458286589419743
$ws.Range("N24").Value = 1.118285595362629
$ws.Range("O24").Value = 2.716275375893474
$ws.Range("B25").Value = 1.350135486690192
$ws.Range("C25").Value = 0.2318216392682757
$ws.Range("D25").Value = 0.02472589985637086
$ws.Range("E25").Value = 0.09818910709614137
$ws.Range("F25").Value = 0.7889064867293314
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("L25").Value = 0.2262667429533423
$ws.Range("N25").Value = 1.136965999057317
$ws.Range("O25").Value = 2.68088948608613
